$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A15").Value = "HSCEI Index"
$ws.Range("A16").Value = "SHSZ300 Index"

$ws.Range("A17").Select()
